$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4810.1
$ws.Range("J40").Value = 4785.7144
$ws.Range("L40").Value = 4785.7144
$ws.Range("N40").Value = -5135.7144

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9870.179
$ws.Range("I74").Value = 9870.179
$ws.Range("K74").Value = 9870.179
$ws.Range("M74").Value = -8934.179

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 9870.179
$ws.Range("I77").Value = 9870.179
$ws.Range("K77").Value = 49350.895
$ws.Range("M77").Value = -44670.895

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 9259758
$ws.Range("I80").Value = 332.26666
$ws.Range("J80").Value = 20834040
$ws.Range("K80").Value = 996.79998
$ws.Range("L80").Value = 62502120
$ws.Range("M80").Value = 1.200019999999995
$ws.Range("N80").Value = -62504116

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 9259758
$ws.Range("I83").Value = 332.26666
$ws.Range("J83").Value = 20834040
$ws.Range("K83").Value = 2990.39994
$ws.Range("L83").Value = 187506360
$ws.Range("M83").Value = 2001.60006
$ws.Range("N83").Value = -187516344

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2590.625
$ws.Range("I98").Value = 1532.2142
$ws.Range("K98").Value = 1532.2142
$ws.Range("M98").Value = -34.21419999999989

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2590.625
$ws.Range("I122").Value = 1532.2142
$ws.Range("K122").Value = 4596.642599999999
$ws.Range("M122").Value = -2146.642599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1265719
$ws.Range("I125").Value = 11363636
$ws.Range("J125").Value = 3479.375
$ws.Range("K125").Value = 102272724
$ws.Range("L125").Value = 31314.375
$ws.Range("M125").Value = -102270264
$ws.Range("N125").Value = -36234.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4926.662
$ws.Range("I138").Value = 7974.5386
$ws.Range("J138").Value = 3275.7292
$ws.Range("K138").Value = 23923.6158
$ws.Range("L138").Value = 9827.187600000001
$ws.Range("M138").Value = -18783.6158
$ws.Range("N138").Value = -20107.1876

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 38373.75
$ws.Range("J37").Value = 64997.5
$ws.Range("L37").Value = 64997.5
$ws.Range("N37").Value = -65543.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1312848.9
$ws.Range("I61").Value = 8934.956
$ws.Range("J61").Value = 2676031.5
$ws.Range("K61").Value = 8934.956
$ws.Range("L61").Value = 2676031.5
$ws.Range("M61").Value = -8722.956
$ws.Range("N61").Value = -2676455.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2308.6667
$ws.Range("I88").Value = 1820
$ws.Range("J88").Value = 2471.5557
$ws.Range("K88").Value = 1820
$ws.Range("L88").Value = 2471.5557
$ws.Range("M88").Value = -1414
$ws.Range("N88").Value = -3283.5557

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2308.6667
$ws.Range("I91").Value = 1820
$ws.Range("J91").Value = 2471.5557
$ws.Range("K91").Value = 1820
$ws.Range("L91").Value = 2471.5557
$ws.Range("M91").Value = -416
$ws.Range("N91").Value = -5279.5557

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 75332.414
$ws.Range("J109").Value = 75332.414
$ws.Range("L109").Value = 75332.414
$ws.Range("N109").Value = -78106.414

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2334.5217
$ws.Range("I132").Value = 2019.9512
$ws.Range("K132").Value = 6059.8536
$ws.Range("M132").Value = -3529.8536

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1312848.9
$ws.Range("I136").Value = 8934.956
$ws.Range("J136").Value = 2676031.5
$ws.Range("K136").Value = 26804.868
$ws.Range("L136").Value = 8028094.5
$ws.Range("M136").Value = -24254.868
$ws.Range("N136").Value = -8033194.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 137871
$ws.Range("J139").Value = 137871
$ws.Range("L139").Value = 137871
$ws.Range("N139").Value = -148151

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4333.6665
$ws.Range("I86").Value = 2125.75
$ws.Range("J86").Value = 8749.5
$ws.Range("K86").Value = 2125.75
$ws.Range("L86").Value = 8749.5
$ws.Range("M86").Value = -1002.75
$ws.Range("N86").Value = -10995.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4333.6665
$ws.Range("I89").Value = 2125.75
$ws.Range("J89").Value = 8749.5
$ws.Range("K89").Value = 10628.75
$ws.Range("L89").Value = 43747.5
$ws.Range("M89").Value = -5012.75
$ws.Range("N89").Value = -54979.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6723.2173
$ws.Range("I99").Value = 8249.117
$ws.Range("K99").Value = 8249.117
$ws.Range("M99").Value = -6751.117

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 15347.7
$ws.Range("I105").Value = 21697.6
$ws.Range("K105").Value = 21697.6
$ws.Range("M105").Value = -19950.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19568420
$ws.Range("I134").Value = 2906
$ws.Range("J134").Value = 81822330
$ws.Range("K134").Value = 8718
$ws.Range("L134").Value = 245466990
$ws.Range("M134").Value = -6183
$ws.Range("N134").Value = -245472060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2214.1428
$ws.Range("I105").Value = 1675
$ws.Range("K105").Value = 1675
$ws.Range("M105").Value = 72

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2415.0625
$ws.Range("I107").Value = 1632.1666
$ws.Range("J107").Value = 4763.75
$ws.Range("K107").Value = 1632.1666
$ws.Range("L107").Value = 4763.75
$ws.Range("M107").Value = 287.8334
$ws.Range("N107").Value = -8603.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9261992
$ws.Range("I132").Value = 2509.0715
$ws.Range("K132").Value = 7527.2145
$ws.Range("M132").Value = -4997.2145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3471.7778
$ws.Range("I109").Value = 974.6667
$ws.Range("K109").Value = 2924.0001
$ws.Range("M109").Value = -1884.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 361.18182
$ws.Range("I2").Value = 123.333336
$ws.Range("J2").Value = 646.6
$ws.Range("K2").Value = 123.333336
$ws.Range("L2").Value = 646.6
$ws.Range("M2").Value = -10.333336
$ws.Range("N2").Value = -872.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9271909
$ws.Range("I80").Value = 191949.95
$ws.Range("K80").Value = 191949.95
$ws.Range("M80").Value = -190951.95

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 9271909
$ws.Range("I83").Value = 191949.95
$ws.Range("K83").Value = 959749.75
$ws.Range("M83").Value = -954757.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5106238
$ws.Range("I132").Value = 1506.9736
$ws.Range("K132").Value = 4520.9208
$ws.Range("M132").Value = -1990.9208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6781.5713
$ws.Range("I82").Value = 992.6667
$ws.Range("J82").Value = 11123.25
$ws.Range("K82").Value = 992.6667
$ws.Range("L82").Value = 11123.25
$ws.Range("M82").Value = -631.6667
$ws.Range("N82").Value = -11845.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 6781.5713
$ws.Range("I85").Value = 992.6667
$ws.Range("J85").Value = 11123.25
$ws.Range("K85").Value = 992.6667
$ws.Range("L85").Value = 11123.25
$ws.Range("M85").Value = 255.3333
$ws.Range("N85").Value = -13619.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 15555.5
$ws.Range("J101").Value = 15555.5
$ws.Range("L101").Value = 15555.5
$ws.Range("N101").Value = -22045.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 47663
$ws.Range("J131").Value = 47663
$ws.Range("L131").Value = 47663
$ws.Range("N131").Value = -57743

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7625.5293
$ws.Range("I132").Value = 5313.92
$ws.Range("J132").Value = 14046.667
$ws.Range("K132").Value = 15941.76
$ws.Range("L132").Value = 42140.001
$ws.Range("M132").Value = -13411.76
$ws.Range("N132").Value = -47200.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4281.6895
$ws.Range("I136").Value = 4118.3477
$ws.Range("K136").Value = 12355.0431
$ws.Range("M136").Value = -9805.043100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 10027.667
$ws.Range("I7").Value = 8973.200000000001
$ws.Range("K7").Value = 8973.200000000001
$ws.Range("M7").Value = -8860.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 25799.6
$ws.Range("I45").Value = 23249
$ws.Range("J45").Value = 27500
$ws.Range("K45").Value = 23249
$ws.Range("L45").Value = 27500
$ws.Range("M45").Value = -22758
$ws.Range("N45").Value = -28482

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 91136
$ws.Range("J92").Value = 91136
$ws.Range("L92").Value = 91136
$ws.Range("N92").Value = -96128

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 83333850
$ws.Range("I107").Value = 361.83334
$ws.Range("K107").Value = 1085.50002
$ws.Range("M107").Value = 834.4999800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3975.2632
$ws.Range("I126").Value = 3934.6924
$ws.Range("K126").Value = 11804.0772
$ws.Range("M126").Value = -9334.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 41051.883
$ws.Range("I132").Value = 64915.625
$ws.Range("J132").Value = 2869.9
$ws.Range("K132").Value = 194746.875
$ws.Range("L132").Value = 8609.700000000001
$ws.Range("M132").Value = -192216.875
$ws.Range("N132").Value = -13669.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5628.303
$ws.Range("I136").Value = 1403.85
$ws.Range("J136").Value = 12127.462
$ws.Range("K136").Value = 4211.549999999999
$ws.Range("L136").Value = 36382.386
$ws.Range("M136").Value = -1661.549999999999
$ws.Range("N136").Value = -41482.386